$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$styleTmp = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.565.87"
$ws.Range("D2").Style = $styleTmp
$ws.Range("E2").Value = "  -6.39%  "

$styleTmp = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.926.29"
$ws.Range("D3").Style = $styleTmp
$ws.Range("E3").Value = "  -8.99%  "

$styleTmp = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = $styleTmp
$ws.Range("E4").Value = "  -0.19%  "

$styleTmp = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.45"
$ws.Range("D5").Style = $styleTmp
$ws.Range("E5").Value = "  -10.28%  "

$styleTmp = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.88"
$ws.Range("D6").Style = $styleTmp
$ws.Range("E6").Value = "  -13.53%  "

$styleTmp = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = $styleTmp
$ws.Range("E7").Value = "  -0.42%  "

$styleTmp = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.916.10"
$ws.Range("D8").Style = $styleTmp
$ws.Range("E8").Value = "  -8.96%  "

$styleTmp = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.457"
$ws.Range("D9").Style = $styleTmp
$ws.Range("E9").Value = "  -16.25%  "

$styleTmp = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("D10").Style = $styleTmp
$ws.Range("E10").Value = "  -19.46%  "

$styleTmp = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.82"
$ws.Range("D11").Style = $styleTmp
$ws.Range("E11").Value = "  -11.53%  "

$styleTmp = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.428"
$ws.Range("D12").Style = $styleTmp
$ws.Range("E12").Value = "  -14.08%  "

$styleTmp = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "31.65"
$ws.Range("D13").Style = $styleTmp
$ws.Range("E13").Value = "  -19.17%  "

$styleTmp = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000199"
$ws.Range("D14").Style = $styleTmp
$ws.Range("E14").Value = "  -18.48%  "

$styleTmp = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.395.46"
$ws.Range("D15").Style = $styleTmp
$ws.Range("E15").Value = "  -9.14%  "

$styleTmp = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.392.64"
$ws.Range("D16").Style = $styleTmp
$ws.Range("E16").Value = "  -6.68%  "

$ws.Range("E17").Value = "  -5.72%  "

$styleTmp = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.904.18"
$ws.Range("D18").Style = $styleTmp
$ws.Range("E18").Value = "  -9.67%  "

$styleTmp = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "467.48"
$ws.Range("D19").Style = $styleTmp
$ws.Range("E19").Value = "  -12.35%  "

$styleTmp = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.13"
$ws.Range("D20").Style = $styleTmp
$ws.Range("E20").Value = "  -14.41%  "

$styleTmp = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.66"
$ws.Range("D21").Style = $styleTmp
$ws.Range("E21").Value = "  -15.71%  "

$styleTmp = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.624"
$ws.Range("D22").Style = $styleTmp
$ws.Range("E22").Value = "  -18.19%  "

$styleTmp = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.38"
$ws.Range("D23").Style = $styleTmp
$ws.Range("E23").Value = "  -19.36%  "

$styleTmp = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.52"
$ws.Range("D24").Style = $styleTmp
$ws.Range("E24").Value = "  -13.10%  "

$styleTmp = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.73"
$ws.Range("D25").Style = $styleTmp
$ws.Range("E25").Value = "  -15.59%  "

$styleTmp = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = $styleTmp
$ws.Range("E26").Value = "  +0.37%  "

$styleTmp = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.49"
$ws.Range("D27").Style = $styleTmp
$ws.Range("E27").Value = "  -22.53%  "

$styleTmp = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.91"
$ws.Range("D28").Style = $styleTmp
$ws.Range("E28").Value = "  -15.61%  "

$styleTmp = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.85"
$ws.Range("D29").Style = $styleTmp
$ws.Range("E29").Value = "  -15.83%  "

$styleTmp = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "24.22"
$ws.Range("D30").Style = $styleTmp
$ws.Range("E30").Value = "  -17.45%  "

$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$styleTmp = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.35"
$ws.Range("D31").Style = $styleTmp
$ws.Range("E31").Value = "  -11.72%  "

$styleTmp = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.04"
$ws.Range("D32").Style = $styleTmp
$ws.Range("E32").Value = "  -9.64%  "

$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$styleTmp = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = $styleTmp
$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$styleTmp = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "465.41"
$ws.Range("D34").Style = $styleTmp
$ws.Range("E34").Value = "  -15.35%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$styleTmp = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.57"
$ws.Range("D35").Style = $styleTmp
$ws.Range("E35").Value = "  -5.46%  "

$styleTmp = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.36"
$ws.Range("D36").Style = $styleTmp
$ws.Range("E36").Value = "  -18.22%  "

$styleTmp = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.66"
$ws.Range("D37").Style = $styleTmp
$ws.Range("E37").Value = "  -18.59%  "

$styleTmp = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0379"
$ws.Range("D38").Style = $styleTmp
$ws.Range("E38").Value = "  -10.94%  "

$styleTmp = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0734"
$ws.Range("D39").Style = $styleTmp
$ws.Range("E39").Value = "  -15.26%  "

$styleTmp = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.111"
$ws.Range("D40").Style = $styleTmp
$ws.Range("E40").Value = "  -11.82%  "

$styleTmp = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.67"
$ws.Range("D41").Style = $styleTmp
$ws.Range("E41").Value = "  -18.22%  "

$styleTmp = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.579.09"
$ws.Range("D42").Style = $styleTmp
$ws.Range("E42").Value = "  -11.64%  "

$ws.Range("E43").Value = "  -0.29%  "

$styleTmp = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.16"
$ws.Range("D44").Style = $styleTmp
$ws.Range("E44").Value = "  -19.47%  "

$styleTmp = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.218"
$ws.Range("D45").Style = $styleTmp
$ws.Range("E45").Value = "  -17.77%  "

$styleTmp = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "110.58"
$ws.Range("D46").Style = $styleTmp
$ws.Range("E46").Value = "  -8.79%  "

$styleTmp = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1000"
$ws.Range("D47").Style = $styleTmp
$ws.Range("E47").Value = "  -12.45%  "

$styleTmp = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.76"
$ws.Range("D48").Style = $styleTmp
$ws.Range("E48").Value = "  -17.51%  "

$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$styleTmp = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0456"
$ws.Range("D49").Style = $styleTmp
$ws.Range("E49").Value = "  -22.35%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$styleTmp = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.27"
$ws.Range("D50").Style = $styleTmp
$ws.Range("E50").Value = "  -19.87%  "

$styleTmp = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.16"
$ws.Range("D51").Style = $styleTmp
$ws.Range("E51").Value = "  -6.87%  "

